$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '27.339.14'
Set-TextValue $ws.Range('E2') '  +1.59%  '
Set-TextValue $ws.Range('D3') '1.833.65'
Set-TextValue $ws.Range('E3') '  +1.06%  '
Set-TextValue $ws.Range('D5') '314.86'
Set-TextValue $ws.Range('E6') '  +0.77%  '
Set-TextValue $ws.Range('D7') '0.4747'
Set-TextValue $ws.Range('E7') '  +2.01%  '
Set-TextValue $ws.Range('D8') '0.3691'
Set-TextValue $ws.Range('E8') '  +0.92%  '
Set-TextValue $ws.Range('D9') '0.07460'
Set-TextValue $ws.Range('E9') '  +1.37%  '
Set-TextValue $ws.Range('D10') '0.8866'
Set-TextValue $ws.Range('E10') '  +2.08%  '
Set-TextValue $ws.Range('D11') '20.46'
Set-TextValue $ws.Range('E11') '  +1.17%  '
Set-TextValue $ws.Range('D12') '1.875.85'
Set-TextValue $ws.Range('E12') '  +2.29%  '
Set-TextValue $ws.Range('D13') '0.07337'
Set-TextValue $ws.Range('E13') '  +3.27%  '
Set-TextValue $ws.Range('D14') '5.440'
Set-TextValue $ws.Range('E14') '  +1.20%  '
Set-TextValue $ws.Range('D15') '93.35'
Set-TextValue $ws.Range('E15') '  +2.64%  '
Set-TextValue $ws.Range('D16') '6.582'
Set-TextValue $ws.Range('E16') '  +1.23%  '
Set-TextValue $ws.Range('E17') '  +0.60%  '
Set-TextValue $ws.Range('D18') '0.000008799'
Set-TextValue $ws.Range('E18') '  +1.17%  '
Set-TextValue $ws.Range('D19') '1.010'
Set-TextValue $ws.Range('E19') '  +0.80%  '
Set-TextValue $ws.Range('B20') 'WrappedBTC'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D20') '27.579.70'
Set-TextValue $ws.Range('E20') '  +2.39%  '
Set-TextValue $ws.Range('B21') 'Avalanche'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D21') '14.81'
Set-TextValue $ws.Range('E21') '  +1.37%  '
Set-TextValue $ws.Range('D22') '5.293'
Set-TextValue $ws.Range('E22') '  +0.04%  '
Set-TextValue $ws.Range('E23') '  +0.90%  '
Set-TextValue $ws.Range('D24') '2.095.00'
Set-TextValue $ws.Range('E24') '  +1.64%  '
Set-TextValue $ws.Range('D25') '1.892'
Set-TextValue $ws.Range('E25') '  -0.15%  '
Set-TextValue $ws.Range('D26') '152.00'
Set-TextValue $ws.Range('E26') '  +0.78%  '
Set-TextValue $ws.Range('D27') '18.64'
Set-TextValue $ws.Range('E27') '  +1.43%  '
Set-TextValue $ws.Range('D28') '2.148'
Set-TextValue $ws.Range('E28') '  +1.03%  '
Set-TextValue $ws.Range('E29') '  +0.08%  '
Set-TextValue $ws.Range('D30') '117.32'
Set-TextValue $ws.Range('E30') '  +1.63%  '
Set-TextValue $ws.Range('D31') '0.08993'
Set-TextValue $ws.Range('E31') '  +0.89%  '
Set-TextValue $ws.Range('D32') '0.7531'
Set-TextValue $ws.Range('E32') '  +0.01%  '
Set-TextValue $ws.Range('D33') '1.177'
Set-TextValue $ws.Range('E33') '  +1.43%  '
Set-TextValue $ws.Range('E34') '  +1.69%  '
Set-TextValue $ws.Range('D35') '2.947'
Set-TextValue $ws.Range('E35') '  +1.31%  '
Set-TextValue $ws.Range('E36') '  +0.85%  '
Set-TextValue $ws.Range('E37') '  +0.66%  '
Set-TextValue $ws.Range('D38') '0.05346'
Set-TextValue $ws.Range('E38') '  +1.30%  '
Set-TextValue $ws.Range('D39') '0.01955'
Set-TextValue $ws.Range('E39') '  +0.71%  '
Set-TextValue $ws.Range('D40') '2.978'
Set-TextValue $ws.Range('E40') '  -0.02%  '
Set-TextValue $ws.Range('D41') '7.277'
Set-TextValue $ws.Range('E41') '  +1.06%  '
Set-TextValue $ws.Range('D42') '2.390'
Set-TextValue $ws.Range('E42') '  +3.92%  '
Set-TextValue $ws.Range('D43') '0.5315'
Set-TextValue $ws.Range('E43') '  +0.63%  '
Set-TextValue $ws.Range('E44') '  +0.47%  '
Set-TextValue $ws.Range('D45') '8.493'
Set-TextValue $ws.Range('E45') '  +0.97%  '
Set-TextValue $ws.Range('E46') '  +1.64%  '
Set-TextValue $ws.Range('D47') '10.51'
Set-TextValue $ws.Range('E47') '  +0.81%  '
Set-TextValue $ws.Range('D48') '105.23'
Set-TextValue $ws.Range('E48') '  +2.04%  '
Set-TextValue $ws.Range('E49') '  +0.88%  '
Set-TextValue $ws.Range('D50') '1.673'
Set-TextValue $ws.Range('E50') '  +0.93%  '
Set-TextValue $ws.Range('D51') '0.06299'
Set-TextValue $ws.Range('E51') '  +0.17%  '
